$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Gehan Adel, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Amira Sobhy'
$ws.Range("G3").Value = 'Dr. Eman Tantawi, Administrator, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Majorelle Magdy'
$ws.Range("G4").Value = 'Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Majorelle Magdy'
$ws.Range("G5").Value = 'Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Asmaa Reda'
$ws.Range("G6").Value = 'Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Mohammad El-Tanany, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy'
$ws.Range("G7").Value = 'Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Kerelos Zareef'
$ws.Range("G11").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G12").Value = 'Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Amira Ibrahim, Dr. Dina Adel, Dr. Marina Youhanna'
$ws.Range("G13").Value = 'Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G20").Value = 'Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range("G25").Value = 'Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil'
$ws.Range("G27").Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range("G30").Value = 'Dr. Shorok Mohammad, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Aya Hanafy'
